$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 1257
$ws.Range("I92").Value = 602.6875
$ws.Range("J92").Value = 3001.8333
$ws.Range("K92").Value = 602.6875
$ws.Range("L92").Value = 3001.8333
$ws.Range("M92").Value = 645.3125
$ws.Range("N92").Value = -5497.8333
$ws.Range("H94").Value = 2414.4443
$ws.Range("I94").Value = 2414.4443
$ws.Range("K94").Value = 2414.4443
$ws.Range("M94").Value = -1963.4443
$ws.Range("H96").Value = 575.3333
$ws.Range("I96").Value = 463
$ws.Range("J96").Value = 800
$ws.Range("K96").Value = 1389
$ws.Range("L96").Value = 2400
$ws.Range("M96").Value = -16
$ws.Range("N96").Value = -5146
$ws.Range("H97").Value = 900
$ws.Range("J97").Value = 900
$ws.Range("L97").Value = 2700
$ws.Range("N97").Value = -3692
$ws.Range("H99").Value = 631
$ws.Range("I99").Value = 407.14285
$ws.Range("J99").Value = 2198
$ws.Range("K99").Value = 1221.42855
$ws.Range("L99").Value = 6594
$ws.Range("M99").Value = 276.5714499999999
$ws.Range("N99").Value = -9590
$ws.Range("H101").Value = 1295.4286
$ws.Range("I101").Value = 490.85715
$ws.Range("J101").Value = 2100
$ws.Range("K101").Value = 1472.57145
$ws.Range("L101").Value = 6300
$ws.Range("M101").Value = 149.4285500000001
$ws.Range("N101").Value = -9544
$ws.Range("H112").Value = 29413160
$ws.Range("I112").Value = 333333730
$ws.Range("J112").Value = 1493.3871
$ws.Range("K112").Value = 1000001190
$ws.Range("L112").Value = 4480.1613
$ws.Range("M112").Value = -1000000082
$ws.Range("N112").Value = -6696.1613
$ws.Range("H129").Value = 808.0513
$ws.Range("I129").Value = 494.26666
$ws.Range("J129").Value = 1004.1667
$ws.Range("K129").Value = 1482.79998
$ws.Range("L129").Value = 3012.5001
$ws.Range("M129").Value = 3517.20002
$ws.Range("N129").Value = -13012.5001
$ws.Range("H132").Value = 1532835
$ws.Range("I132").Value = 1490.3793
$ws.Range("K132").Value = 4471.1379
$ws.Range("M132").Value = -1941.1379
$ws.Range("H138").Value = 2780411.5
$ws.Range("I138").Value = 2179
$ws.Range("J138").Value = 3706488.8
$ws.Range("K138").Value = 6537
$ws.Range("L138").Value = 11119466.4
$ws.Range("M138").Value = -1397
$ws.Range("N138").Value = -11129746.4

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 47716270
$ws.Range("I61").Value = 100101330
$ws.Range("J61").Value = 93492.63
$ws.Range("K61").Value = 100101330
$ws.Range("L61").Value = 93492.63
$ws.Range("M61").Value = -100101118
$ws.Range("N61").Value = -93916.63
$ws.Range("H74").Value = 7638044.5
$ws.Range("J74").Value = 113611.11
$ws.Range("L74").Value = 113611.11
$ws.Range("N74").Value = -115359.11
$ws.Range("H77").Value = 7638044.5
$ws.Range("J77").Value = 113611.11
$ws.Range("L77").Value = 568055.55
$ws.Range("N77").Value = -576791.55
$ws.Range("H136").Value = 47716270
$ws.Range("I136").Value = 100101330
$ws.Range("J136").Value = 93492.63
$ws.Range("K136").Value = 300303990
$ws.Range("L136").Value = 280477.89
$ws.Range("M136").Value = -300301440
$ws.Range("N136").Value = -285577.89

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1033.8077
$ws.Range("I99").Value = 1073.95
$ws.Range("J99").Value = 900
$ws.Range("K99").Value = 1073.95
$ws.Range("L99").Value = 900
$ws.Range("M99").Value = 424.05
$ws.Range("N99").Value = -3896

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 62503572
$ws.Range("I58").Value = 76926270
$ws.Range("J58").Value = 5200.3335
$ws.Range("K58").Value = 76926270
$ws.Range("L58").Value = 5200.3335
$ws.Range("M58").Value = -76926067
$ws.Range("N58").Value = -5606.3335
$ws.Range("H132").Value = 27244.75
$ws.Range("I132").Value = 1946.9667
$ws.Range("J132").Value = 103138.1
$ws.Range("K132").Value = 5840.9001
$ws.Range("L132").Value = 309414.3
$ws.Range("M132").Value = -3310.9001
$ws.Range("N132").Value = -314474.3
$ws.Range("H136").Value = 62503572
$ws.Range("I136").Value = 76926270
$ws.Range("J136").Value = 5200.3335
$ws.Range("K136").Value = 230778810
$ws.Range("L136").Value = 15601.0005
$ws.Range("M136").Value = -230776260
$ws.Range("N136").Value = -20701.0005
$ws.Range("H141").Value = 49132.273
$ws.Range("J141").Value = 49132.273
$ws.Range("L141").Value = 49132.273
$ws.Range("N141").Value = -59492.273

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 850.0741
$ws.Range("I5").Value = 840
$ws.Range("J5").Value = 852.36365
$ws.Range("K5").Value = 2520
$ws.Range("L5").Value = 2557.09095
$ws.Range("M5").Value = -2408
$ws.Range("N5").Value = -2781.09095
$ws.Range("H68").Value = 488875.66
$ws.Range("I68").Value = 979.37933
$ws.Range("J68").Value = 1667958.4
$ws.Range("K68").Value = 2938.13799
$ws.Range("L68").Value = 5003875.199999999
$ws.Range("M68").Value = -2127.13799
$ws.Range("N68").Value = -5005497.199999999
$ws.Range("H71").Value = 488875.66
$ws.Range("I71").Value = 979.37933
$ws.Range("J71").Value = 1667958.4
$ws.Range("K71").Value = 8814.41397
$ws.Range("L71").Value = 15011625.6
$ws.Range("M71").Value = -4758.41397
$ws.Range("N71").Value = -15019737.6
$ws.Range("H107").Value = 440.71112
$ws.Range("I107").Value = 724.2
$ws.Range("J107").Value = 298.96667
$ws.Range("K107").Value = 2172.6
$ws.Range("L107").Value = 896.9000100000001
$ws.Range("M107").Value = -252.6000000000004
$ws.Range("N107").Value = -4736.90001
$ws.Range("H122").Value = 1000.7857
$ws.Range("J122").Value = 1663.1666
$ws.Range("L122").Value = 14968.4994
$ws.Range("N122").Value = -19868.4994
$ws.Range("H135").Value = 850.0741
$ws.Range("I135").Value = 840
$ws.Range("J135").Value = 852.36365
$ws.Range("K135").Value = 7560
$ws.Range("L135").Value = 7671.27285
$ws.Range("M135").Value = -5025
$ws.Range("N135").Value = -12741.27285
$ws.Range("H137").Value = 18923.295
$ws.Range("J137").Value = 28726.092
$ws.Range("L137").Value = 86178.276
$ws.Range("N137").Value = -96378.276
$ws.Range("H138").Value = 3202.3333
$ws.Range("I138").Value = 2561
$ws.Range("J138").Value = 3785.3635
$ws.Range("K138").Value = 7683
$ws.Range("L138").Value = 11356.0905
$ws.Range("M138").Value = -2543
$ws.Range("N138").Value = -21636.0905

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1237.2727
$ws.Range("I97").Value = 1401.25
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 1401.25
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -905.25
$ws.Range("N97").Value = -1792

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 63342.438
$ws.Range("I100").Value = 45879.637
$ws.Range("K100").Value = 91759.274
$ws.Range("M100").Value = -91218.274
